# Edit script: thêm đồng hồ vào csdl, hiển thị đồng hồ độc quyền, phái đẹp,
# phái mạnh, trẻ em ra màn hình user
#
# Strategy:
#  1. Locate the paragraph containing "MVW: 1" (last brand under "độc quyền").
#  2. Change its trailing "X: 1" run into "X: " + "3" (two runs) -> "MVW: 3".
#  3. Insert 9 new list paragraphs right after it, all sharing the same
#     ListParagraph style + numId=4 numbering used throughout this list:
#       - DAUMIER: 1          (ilvl 7)
#       - ELIO: 1             (ilvl 7)
#       - KORLEX : 1          (ilvl 7)
#       - Phái đẹp:           (ilvl 6, new sub-heading)
#       - Nakzen: 7           (ilvl 7)
#       - Phái mạnh:          (ilvl 6, new sub-heading)
#       - Nakzen: 3           (ilvl 7)
#       - trẻ em              (ilvl 6, new sub-heading)
#       - SKMEI: 6            (ilvl 7)

$d = $word.ActiveDocument

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Step 1: find the "MVW: 1" paragraph -----------------------------------
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd() -eq "MVW: 1") {
        $target = $cand
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'MVW: 1' paragraph"
}

# --- Step 2: rewrite that paragraph's content so it reads "MVW: 3" ---------
# (keeps the same pPr / numbering level, just edits the run content so the
#  trailing ": 1" run becomes ": " + a new "3" run)
$rewrittenPara = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="7"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>MV</w:t></w:r><w:r><w:t>W</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>3</w:t></w:r></w:p>'
$target.Range.InsertXML($xmlHeader + $rewrittenPara + $xmlFooter)

# Re-resolve the (now rewritten) paragraph so we can insert right after it.
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd() -eq "MVW: 3") {
        $target = $cand
        break
    }
}
if ($target -eq $null) {
    throw "Could not locate the rewritten 'MVW: 3' paragraph"
}

# --- Step 3: insert the 9 new paragraphs after it ---------------------------
# Create one placeholder paragraph, then swap its content for all 9 new
# paragraphs at once via InsertXML (InsertXML replaces the whole paragraph
# that the (now collapsed) range belongs to, so a fresh empty paragraph is
# the safe place to target).
$targetIndex = $target.Index
$target.Range.InsertParagraphAfter()

$placeholder = $d.Paragraphs.Item($targetIndex + 1)
if ($placeholder.Range.Text.TrimEnd() -ne "") {
    throw "Placeholder paragraph at index $($targetIndex + 1) was not empty"
}

$newParas = ""
$newParas += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="7"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>DAUMIER</w:t></w:r><w:r><w:t>: 1</w:t></w:r></w:p>'
$newParas += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="7"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>ELIO: 1</w:t></w:r></w:p>'
$newParas += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="7"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">KORLEX </w:t></w:r><w:r><w:t>: 1</w:t></w:r></w:p>'
$newParas += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="6"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Ph</w:t></w:r><w:r><w:t>ái</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>đẹp</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>'
$newParas += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="7"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Nakzen: </w:t></w:r><w:r><w:t>7</w:t></w:r></w:p>'
$newParas += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="6"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Ph</w:t></w:r><w:r><w:t>ái</w:t></w:r><w:r><w:t xml:space="preserve"> m</w:t></w:r><w:r><w:t>ạnh</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>'
$newParas += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="7"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Nakzen</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>3</w:t></w:r></w:p>'
$newParas += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="6"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>tr</w:t></w:r><w:r><w:t>ẻ</w:t></w:r><w:r><w:t xml:space="preserve"> em</w:t></w:r></w:p>'
$newParas += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="7"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>SKMEI</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>6</w:t></w:r></w:p>'

$placeholder.Range.InsertXML($xmlHeader + $newParas + $xmlFooter)

Write-Output "Done. Paragraph count now: $($d.Paragraphs.Count)"
